$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row correct-answer marks value (B11: 3 -> 5)
$ws.Range("B11").Value = 5

# Update "Total" row correct-answer marks value (B12: 75 -> 125)
$ws.Range("B12").Value = 125

# Update total score display (E12: "73/84" -> "125/140")
$ws.Range("E12").Value = "125/140"
